$wb = $excel.ActiveWorkbook

# "Metadata" sheet: update URL and Date values
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://nphcda.gov.ng/ig/fhir/ValueSet/ng-gender"
$wsMeta.Range("B8").Value = "2025-07-03T11:54:34+01:00"

# "Include #0" sheet: update System URI value (row 4, column B)
$wsInclude = $wb.Worksheets.Item("Include #0")
$wsInclude.Range("B4").Value = "https://nphcda.gov.ng/immunizationIG/CodeSystem/nigeria-gender"
